$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the three oldest year rows (2007年, 2008年, 2009年) which are
# currently in rows 2-4. Remaining rows (2010年-2013年) shift up to
# become rows 2-5.
$ws.Range("A2:F4").Delete()
